# "Generate Report for Handback"
# The 569d02c7-ed69-4da3-bcea-4a677ba8dd86.md file has just been handed back
# (on both the zh-cn and de-de locales), so its status flips from
# "Ready for handoff" to "Handed back: in sync with en-US" and the
# "Latest Handback DateTime" column is stamped with the handback time.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------
# Row 3 corresponds to 569d02c7-ed69-4da3-bcea-4a677ba8dd86.md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $statusHandedBack
$wsZhCn.Range("G3").Value = "2016-03-03 07:41:48"

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $statusHandedBack
$wsDeDe.Range("G3").Value = "2016-03-03 07:42:10"
